# Generate Report for Handback
# The b6e66154 file has now been handed back (it was previously "Ready for
# handoff"), and the "Latest Handback DateTime" for both files is refreshed
# to reflect the new handback timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the b6e66154 file, status now "handed back" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("G2").Value = "2016-03-09 10:22:29"
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-03-09 10:22:29"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("G2").Value = "2016-03-09 10:22:47"
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-03-09 10:22:47"
